$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage on the numeric-looking columns so they stay strings
# instead of being auto-converted to numbers by Excel
$ws.Range("G2:K2").NumberFormat = "@"

# Copy the row 5 data (Nov 1 2020 match) into row 2
$ws.Range("A2").Value = " Nov 1 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "KL Rahul (c)†"
$ws.Range("G2").Value = "29"
$ws.Range("H2").Value = "27"
$ws.Range("I2").Value = "3"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "107.40"

# Delete the now-redundant rows 3 through 6
$ws.Range("A3:K6").EntireRow.Delete()
